$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (P) mirroring the existing "2020" column (O):
# copy O4/O5 formatting into P4/P5, then overwrite with the new values so
# the new cells keep the same number format / font / borders / alignment
# as the rest of the year header row and data row.

$ws.Range("O4").Copy($ws.Range("P4"))
$ws.Range("P4").Value = 2021

$ws.Range("O5").Copy($ws.Range("P5"))
$ws.Range("P5").Value = 80.9

# Move the active selection (as recorded in the saved view state).
$ws.Range("N10").Select()
